$d = $word.ActiveDocument

$replacements = @(
    @("927×5=4635", "294×9=2646"),
    @("942×9=8478", "446×4=1784"),
    @("178×6=1068", "805×5=4025"),
    @("949×7=6643", "739×2=1478"),
    @("146×8=1168", "736×3=2208"),
    @("648×6=3888", "631×4=2524"),
    @("114×9=1026", "809×4=3236"),
    @("696×7=4872", "578×4=2312"),
    @("817×4=3268", "831×3=2493"),
    @("518×7=3626", "335×6=2010"),
    @("406×7=2842", "993×5=4965"),
    @("221×6=1326", "464×9=4176"),
    @("818×4=3272", "947×4=3788"),
    @("849×9=7641", "855×6=5130"),
    @("111×2=222",  "724×2=1448"),
    @("545×9=4905", "404×4=1616"),
    @("265×7=1855", "525×3=1575"),
    @("975×6=5850", "865×5=4325"),
    @("255×2=510",  "283×4=1132"),
    @("734×5=3670", "900×8=7200"),
    @("811×8=6488", "153×2=306"),
    @("400×5=2000", "287×8=2296"),
    @("905×4=3620", "435×5=2175"),
    @("589×5=2945", "364×4=1456"),
    @("566×7=3962", "691×6=4146")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "done"
